# Updated cryptos list on Thu Jun 29 14:37:13 UTC 2023 with GitHub Actions
# Refreshes Price (D) / Volume(1h) (E) figures, and swaps the row order of a
# couple of coin pairs (B/C/D/E) that changed rank between runs.
#
# Note: several "Price" values look numeric (e.g. "1.001", "1.0000",
# "0.000007361") but must stay as literal text (matching the source feed's
# formatting, including trailing zeros). Setting .NumberFormat = "@" before
# assigning .Value keeps Excel from coercing them into numbers, and resetting
# .Style back to "Normal" afterwards avoids leaving a stray number-format
# style on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.587.42'
$ws.Range("E2").Value = '  +0.96%  '

$ws.Range("D3").Value = '1.858.87'
$ws.Range("E3").Value = '  -0.09%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.29%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4704'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.40%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2748'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.22%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06351'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.23%  '

$ws.Range("E10").Value = '  +8.07%  '

$ws.Range("D11").Value = '1.848.28'
$ws.Range("E11").Value = '  -0.39%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07442'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.150'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.80%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '84.94'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.76%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6309'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.49%  '

$ws.Range("D16").Value = '30.577.94'
$ws.Range("E16").Value = '  +1.13%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '243.12'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.0000'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.81'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.41%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007361'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.28%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.21%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.990'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.024'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.326'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.61%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '164.90'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.90%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.889'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1016'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.51%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.381'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.19%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.057'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.11%  '

$ws.Range("E31").Value = '  -1.80%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.04914'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.33%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.149'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7066'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.58%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.710'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.72%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.01915'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.36%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.689'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.94%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.8797'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.88%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.991'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.84%  '

$ws.Range("E40").Value = '  -0.37%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.001'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.24%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4086'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.81%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.539'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.21%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.284'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.43%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.35'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.57%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1215'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.20%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.643'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.76%  '

$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '33.51'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.58%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05547'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.371'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.59%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3700'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.47%  '
